$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Introduce new shared strings in the exact order required so the ---
# --- rebuilt sharedStrings.xml table lands in the desired sequence:   ---
# --- BOOLEAN, DATE, STRING, "testing    ", INT, FLOAT, STRDATE, 2022-01-01

# Row 6 (STRING / "testing    ") gets written first so those two new
# strings are appended before the others.
$ws.Range("A6").Value = "STRING"
$ws.Range("B6").Value = "testing    "

# Row 2 becomes INT (replaces the old "NUMBER" label).
$ws.Range("A2").Value = "INT"

# Row 3 becomes FLOAT (replaces the old "DATE" label that used to live here).
$ws.Range("A3").Value = "FLOAT"

# Row 5 becomes STRDATE, with its value stored as literal text "2022-01-01".
$ws.Range("A5").Value = "STRDATE"
$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = "2022-01-01"

# Row 4 becomes the relocated DATE row (label only for now).
$ws.Range("A4").Value = "DATE"

# --- Now fill in the numeric values / number formats for column B ---

# B2: INT value, formatted with thousands separator.
$ws.Range("B2").NumberFormat = "#,##0"
$ws.Range("B2").Value = 10100

# B3: FLOAT value, plain General formatting (clear the inherited date format).
$ws.Range("B3").ClearFormats()
$ws.Range("B3").Value = 5.5

# B4: DATE value (moved down from the old row 3), keep the custom date format.
$ws.Range("B4").NumberFormat = "mm/dd/yy;@"
$ws.Range("B4").Value = 42615

# --- Restore the selected cell shown in the saved workbook ---
$ws.Range("D15").Select() | Out-Null
